$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")

# Insert a new row right before the existing "lock" / "no" row (row 8),
# pushing it (and the trailing blank row) down by one.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row with the new "style" / "default" setting.
$ws.Cells.Item(8, 1).Value = "style"
$ws.Cells.Item(8, 2).Value = "default"
